$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row below the existing data row (row 16), which pushes
#    the blank spacer rows and the signature-block rows (old 21/22) down
#    by one (new 22/23), matching the diff.
$ws.Rows.Item(17).Insert()

# 2) Clone the formatting of row 16 (the existing worker/period line) onto
#    the newly inserted row 17 so borders/fonts/number formats match.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new period row (same worker, new period 2509).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45552629"
$ws.Range("D17").Value = "DARLYS PATERNINA CANTILLO"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# 4) Update the summary figures: total overdue value and period count.
$ws.Range("E11").Value = 104390
$ws.Range("F13").Value = 2
